# Update the "Corr/total marks" figures on the marksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Marking row: correct answer marks value 3 -> 5
$ws.Range("B11").Value = 5

# Total row: total marks 45 -> 75, and out-of total 45/84 -> 75/140
$ws.Range("B12").Value = 75
$ws.Range("E12").Value = "75/140"

$wb.Save()
